$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts old rows 13-21 down to 14-22, carrying
# their row heights and cell styles along automatically).
$ws.Rows.Item(13).Insert()

# The insert duplicated row 12's column-A formatting into the new A13; the
# target layout has no content/formatting there at all, so clear it.
$ws.Cells.Item(13, 1).Clear()

# New B13/C13 cells don't exist yet, so copy the (wrap/vertical-top) number
# formats from an existing B/C data cell first, then set the value - this
# reuses the workbook's existing style indices instead of minting new ones.
$ws.Cells.Item(9, 2).Copy()
$ws.Cells.Item(13, 2).PasteSpecial(-4122)
$ws.Cells.Item(9, 3).Copy()
$ws.Cells.Item(13, 3).PasteSpecial(-4122)

# --- Update Objetivos: (row 10) body text, column B and C ---
$objetivos = "Desenvolver um projeto sobre tema de Engenharia de Produção, similar a situações que os alunos irão encontrar na vida real no efetivo exercício de sua profissão, `nAplicar e integrar conhecimentos adquiridos em demais disciplinas de seu curso`nDesenvolver competências técnicas, as relacionadas ao projeto em si, bem como competências transversais (habilidades e atitudes), num ambiente de aprendizagem baseado em PBL (Project-Baed Learning)."
$ws.Range("B10:C10").Value = $objetivos

# --- New row 13: professor name under "Docentes responsáveis:" ---
$professor = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("B13:C13").Value = $professor

# --- Programa resumido: (row 14) ---
$programaResumido = "Tópicos que abordem o tema do projeto de seu planejamento a execução."
$ws.Range("B14:C14").Value = $programaResumido

# --- Programa: (row 16) ---
$programa = "Noções de Gestão de Projetos`nOrganização do tempo: dimensão pessoal;`nTécnicas para a realização de apresentações;`nNoções de Aprendizagem Baseada em Projetos`nTrabalho em Grupo, Equipes e times. `nPostura e Ética Profissional`nTécnicas para redação de relatório técnico;`nTutoria de projetos.`nAssuntos Técnicos específicos relacionados com o tema do projeto."
$ws.Range("B16:C16").Value = $programa

# --- Método: (row 19) ---
$metodo = "O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras.`n`nOs alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Produção, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão. `nCada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.`nAs aulas ocorrerão: 1) através de uma reunião da equipe de trabalho para tratar do projeto, e  2) palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores  ou profissionais de empresas."
$ws.Range("B19:C19").Value = $metodo

# --- Critério: (row 20) ---
$criterio = "A nota será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.`nO detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na coordenação da disciplina."
$ws.Range("B20:C20").Value = $criterio

# --- Norma de recuperação: (row 21) ---
$norma = "Não há recuperação"
$ws.Range("B21:C21").Value = $norma

# --- Bibliografia: (row 22) ---
$bibliografia = "Artigos sobre metodologias ativas de aprendizagem e  Project Based Learning.`nLivros e Artigos científicos relacionados com o tema do projeto."
$ws.Range("B22:C22").Value = $bibliografia

Write-Output "done"
